$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.976.52'
$ws.Range("E2").Value = '  -1.00%  '

$ws.Range("D3").Value = '2.005.57'
$ws.Range("E3").Value = '  -2.33%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.993'
$ws.Range("E4").Value = '  -1.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.90'
$ws.Range("E5").Value = '  -2.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.605'
$ws.Range("E6").Value = '  -2.62%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.45'
$ws.Range("E8").Value = '  -4.65%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.377'
$ws.Range("E9").Value = '  -1.74%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0784'
$ws.Range("E10").Value = '  +1.57%  '

$ws.Range("E11").Value = '  -2.88%  '

$ws.Range("D12").Value = '2.307.08'
$ws.Range("E12").Value = '  -1.99%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.12'
$ws.Range("E13").Value = '  -3.60%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.16'
$ws.Range("E14").Value = '  -2.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.738'
$ws.Range("E15").Value = '  -2.66%  '

$ws.Range("E16").Value = '  -3.15%  '

$ws.Range("D17").Value = '1.994.76'
$ws.Range("E17").Value = '  -2.65%  '

$ws.Range("D18").Value = '36.877.54'
$ws.Range("E18").Value = '  -1.22%  '

$ws.Range("E19").Value = '  +0.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.59'
$ws.Range("E20").Value = '  -1.64%  '

$ws.Range("D21").Value = '0.0₃0814'
$ws.Range("E21").Value = '  -1.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '222.97'
$ws.Range("E22").Value = '  -1.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.14%  '

$ws.Range("E24").Value = '  +1.65%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.18'
$ws.Range("E25").Value = '  -6.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.42'
$ws.Range("E26").Value = '  -2.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.13'
$ws.Range("E27").Value = '  -6.74%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.65'
$ws.Range("E28").Value = '  -2.85%  '

$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.124'
$ws.Range("E29").Value = '  -4.35%  '

$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.34'
$ws.Range("E30").Value = '  -0.25%  '

$ws.Range("E31").Value = '  -3.96%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.49'
$ws.Range("E32").Value = '  -0.87%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0609'
$ws.Range("E33").Value = '  -2.19%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.41'
$ws.Range("E34").Value = '  -3.46%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.33'
$ws.Range("E35").Value = '  -6.01%  '

$ws.Range("E36").Value = '  +1.79%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.32%  '

$ws.Range("E38").Value = '  -4.47%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.35'
$ws.Range("E39").Value = '  +0.70%  '

$ws.Range("D40").Value = '1.470.71'
$ws.Range("E40").Value = '  -0.78%  '

$ws.Range("E41").Value = '  -4.50%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '16.50'
$ws.Range("E42").Value = '  -0.40%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '94.47'
$ws.Range("E43").Value = '  -3.87%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0918'
$ws.Range("E44").Value = '  -3.57%  '

$ws.Range("E45").Value = '  -5.19%  '

$ws.Range("E46").Value = '  -4.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.18'
$ws.Range("E47").Value = '  -0.86%  '

$ws.Range("E48").Value = '  -2.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.90'
$ws.Range("E49").Value = '  -1.48%  '

$ws.Range("D50").Value = '2.195.68'
$ws.Range("E50").Value = '  -2.08%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.21'
$ws.Range("E51").Value = '  -2.72%  '
